# Endtest report template: add a "Result" column between the "Data Matrix"
# and "Created By" columns, and tidy up the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F - this shifts the existing F:J content
# (Data Matrix-value / Created By / Date Created / helper columns) one to
# the right (becoming G:K) and extends the B2:G2 merged banner to B2:H2
# automatically.
$ws.Columns("F:F").Insert()

# Match the new column's width to the template's narrower helper columns.
$ws.Columns("F:F").ColumnWidth = 15.33

# New header for the inserted column, styled like its header-row neighbours.
$ws.Range("F8").Value = "Result"

# Leave the selection where the author left it when saving.
$ws.Range("B11").Select()
